# NPC.xlsx edit: "add NPC move type"
# Renames the old "DropPack" column header to "DropPackList", and inserts
# two new columns (MoveType, AtkDis) before it, with data for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at W:X - this pushes the existing
# "DropPack" column (currently W) to Y, preserving its content/shared string.
$ws.Columns("W:X").Insert()

# Rename the (now relocated) DropPack header to DropPackList first, while it
# is still the sole reference to shared string "DropPack" - this keeps the
# rename in place (same shared-string slot) instead of minting a new one.
$ws.Range("Y1").Value = "DropPackList"

# New column headers
$ws.Range("W1").Value = "MoveType"
$ws.Range("X1").Value = "AtkDis"

# The "AtkDis" header is rich text in the target: the leading "A" keeps the
# default run formatting while "tkDis" carries an explicit font run (11pt,
# black, 宋体, family 3) - reproduce that with a Characters() sub-range.
$xChars = $ws.Range("X1").Characters(2, 5)
$xChars.Font.Name = "宋体"
$xChars.Font.Size = 11
$xChars.Font.ColorIndex = 1
$xChars.Font.Family = 3

# The cell itself also picks up the new style (fontId 2 in the target
# cellXfs) for its default run / cell-level formatting.
$xFont = $ws.Range("X1").Font
$xFont.Name = "宋体"
$xFont.Size = 11
$xFont.ColorIndex = 1
$xFont.Family = 3

# Fill in the MoveType / AtkDis values for the 5 data rows.
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 20

$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 20

$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 20

$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 20

$ws.Range("W6").Value = 2
$ws.Range("X6").Value = 20

# Column widths: V shrinks from 32.125 to 25, W/X get width 25, Y gets ~13.875.
$ws.Columns("V").ColumnWidth = 24.29
$ws.Columns("W").ColumnWidth = 24.29
$ws.Columns("X").ColumnWidth = 24.29
$ws.Columns("Y").ColumnWidth = 13.15

# Move the view/selection the way the author left it.
$ws.Range("X10").Select()
